{"js": "// Apply the three wording fixes in the \"Th\u1ed1ng k\u00ea doanh thu\" use case table\n// (post peer-review edits):\n//   1. \"ti\u1ec1n c\u1ecdc NPP \u0111\u00f3ng, ti\u1ec1n\" -> \"ti\u1ec1n NPP tr\u1ea3 khi nh\u1eadn h\u00e0ng, ti\u1ec1n\"\n//   2. \"h\u00f3a \u0111\u01a1n c\u00f4ng n\u1ee3, phi\u1ebfu c\u00f4ng n\u1ee3 gi\u1eefa\" -> \"h\u00f3a \u0111\u01a1n c\u00f4ng n\u1ee3, phi\u1ebfu chi gi\u1eefa\"\n//   3. \"Doanh thu = ti\u1ec1n c\u1ecdc + ti\u1ec1n\" -> \"Doanh thu = ti\u1ec1n giao h\u00e0ng thu \u0111\u01b0\u1ee3c + ti\u1ec1n\"\n\nconst edits = [\n  {\n    find: \"ti\u1ec1n c\u1ecdc NPP \u0111\u00f3ng, ti\u1ec1n\",\n    replace: \"ti\u1ec1n NPP tr\u1ea3 khi nh\u1eadn h\u00e0ng, ti\u1ec1n\",\n  },\n  {\n    find: \"h\u00f3a \u0111\u01a1n c\u00f4ng n\u1ee3, phi\u1ebfu c\u00f4ng n\u1ee3 gi\u1eefa\",\n    replace: \"h\u00f3a \u0111\u01a1n c\u00f4ng n\u1ee3, phi\u1ebfu chi gi\u1eefa\",\n  },\n  {\n    find: \"Doanh thu = ti\u1ec1n c\u1ecdc + ti\u1ec1n\",\n    replace: \"Doanh thu = ti\u1ec1n giao h\u00e0ng thu \u0111\u01b0\u1ee3c + ti\u1ec1n\",\n  },\n];\n\nfor (const { find, replace } of edits) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply the three wording fixes in the \"Th\u1ed1ng k\u00ea doanh thu\" use case table\n# (post peer-review edits):\n#   1. \"ti\u1ec1n c\u1ecdc NPP \u0111\u00f3ng, ti\u1ec1n\" -> \"ti\u1ec1n NPP tr\u1ea3 khi nh\u1eadn h\u00e0ng, ti\u1ec1n\"\n#   2. \"h\u00f3a \u0111\u01a1n c\u00f4ng n\u1ee3, phi\u1ebfu c\u00f4ng n\u1ee3 gi\u1eefa\" -> \"h\u00f3a \u0111\u01a1n c\u00f4ng n\u1ee3, phi\u1ebfu chi gi\u1eefa\"\n#   3. \"Doanh thu = ti\u1ec1n c\u1ecdc + ti\u1ec1n\" -> \"Doanh thu = ti\u1ec1n giao h\u00e0ng thu \u0111\u01b0\u1ee3c + ti\u1ec1n\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1              # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n\n    $ok = $find.Execute(\n        $findText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $replaceText,\n        2                        # wdReplaceAll\n    )\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $findText\"\n    }\n}\n\nReplace-Text \"ti\u1ec1n c\u1ecdc NPP \u0111\u00f3ng, ti\u1ec1n\" \"ti\u1ec1n NPP tr\u1ea3 khi nh\u1eadn h\u00e0ng, ti\u1ec1n\"\nReplace-Text \"h\u00f3a \u0111\u01a1n c\u00f4ng n\u1ee3, phi\u1ebfu c\u00f4ng n\u1ee3 gi\u1eefa\" \"h\u00f3a \u0111\u01a1n c\u00f4ng n\u1ee3, phi\u1ebfu chi gi\u1eefa\"\nReplace-Text \"Doanh thu = ti\u1ec1n c\u1ecdc + ti\u1ec1n\" \"Doanh thu = ti\u1ec1n giao h\u00e0ng thu \u0111\u01b0\u1ee3c + ti\u1ec1n\"\n"}
